$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(33, 40, 41, 45, 48, 50, 54, 74, 81, 82, 84, 85, 87, 92, 106, 117, 131, 133, 161, 172, 173, 176, 177, 178, 180, 183, 184, 187, 204, 208, 210, 226, 229, 238, 244, 263, 265, 291, 313, 316, 339, 376, 383, 388, 411, 431, 447, 451, 457, 469, 473, 475, 484)
$damslTag = @('sv', 'aa', 'b', 'sd', 'sv', 'aa', 'b', 'sd', 'sv', 'sd', 'ba', 'aa', 'sv', 'sd', 'ba', 'sv', 'sd', 'b', 'sv', 'sd', 'sd', 'sd', 'sd', 'sv', 'sv', '%', 'sd', 'sv', 'sv', 'b', 'sv', '%', 'b', 'sv', 'b', 'ba', 'ba', 'ba', 'b', 'b', 'sd', 'sv', 'b', 'ba', 'b', 'sd', 'sv', 'sv', 'ba', 'b', 'b', 'sd', 'sd')
$dialogAct = @('Statement-opinion', 'Agree/Accept', 'Acknowledge (Backchannel)', 'Statement-non-opinion', 'Statement-opinion', 'Agree/Accept', 'Acknowledge (Backchannel)', 'Statement-non-opinion', 'Statement-opinion', 'Statement-non-opinion', 'Appreciation', 'Agree/Accept', 'Statement-opinion', 'Statement-non-opinion', 'Appreciation', 'Statement-opinion', 'Statement-non-opinion', 'Acknowledge (Backchannel)', 'Statement-opinion', 'Statement-non-opinion', 'Statement-non-opinion', 'Statement-non-opinion', 'Statement-non-opinion', 'Statement-opinion', 'Statement-opinion', 'Uninterpretable', 'Statement-non-opinion', 'Statement-opinion', 'Statement-opinion', 'Acknowledge (Backchannel)', 'Statement-opinion', 'Uninterpretable', 'Acknowledge (Backchannel)', 'Statement-opinion', 'Acknowledge (Backchannel)', 'Appreciation', 'Appreciation', 'Appreciation', 'Acknowledge (Backchannel)', 'Acknowledge (Backchannel)', 'Statement-non-opinion', 'Statement-opinion', 'Acknowledge (Backchannel)', 'Appreciation', 'Acknowledge (Backchannel)', 'Statement-non-opinion', 'Statement-opinion', 'Statement-opinion', 'Appreciation', 'Acknowledge (Backchannel)', 'Acknowledge (Backchannel)', 'Statement-non-opinion', 'Statement-non-opinion')

for ($k = 0; $k -lt $rows.Count; $k++) {
    $r = $rows[$k]
    $ws.Cells.Item($r, 9).Value = $damslTag[$k]
    $ws.Cells.Item($r, 10).Value = $dialogAct[$k]
}
